$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column A (the old "input_useremail" column) and shift column B
# ("input_useremail_1") left so it becomes the sole remaining column.
$ws.Columns("A").Delete()

# Now set the surviving column's header text/value and width to match
# the target state ("input_platformUserAutocomplete", width 32).
$ws.Range("A1").Value = "input_platformUserAutocomplete"
$ws.Columns("A").ColumnWidth = 31.16666666666667
